# Update countries & provincias Spain
# - Swap ranking positions of "Corea del Sur" and "Bosnia y Herzegovina"
#   (row 76 becomes Bosnia y Herzegovina, row 77 becomes Corea del Sur)
# - Refresh the daily COVID stats (Casos totales, Nuevos casos, Casos
#   activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#   rows whose figures moved
# - Bump the "Datos actualizados" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country names at rows 76/77 swap order (Bosnia now ranks above Corea del Sur)
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("A77").Value = "Corea del Sur"

# row -> @{ col = value }
$updates = @{
    4   = @{ B = 6176105; C = 2869;  D = 3425925; E = 2562944;              G = 12;  H = 187236 }
    6   = @{ B = 3641048; C = 21879; D = 2791264; E = 784977;               G = 190; H = 64807 }
    14  = @{                        D = 301195;  E = 98733;                G = 41;  H = 8498 }
    23  = @{ B = 243774;  C = 479;               E = 16566;                G = 2;   H = 9366 }
    29  = @{ B = 118778;  C = 203;   D = 115667;  E = 2914 }
    42  = @{ B = 84379;                                                    G = 5;   H = 5808 }
    51  = @{ B = 58012;   C = 244;   D = 41961;   E = 14229;               G = 3;   H = 1822 }
    62  = @{ B = 41787;   C = 363;   D = 39068;   E = 2400;                G = 6;   H = 319 }
    66  = @{ B = 36920;   C = 220;   D = 25872;   E = 10053;               G = 3;   H = 995 }
    67  = @{ B = 36435;   C = 126;   D = 33847;   E = 2054;                G = 3;   H = 534 }
    69  = @{ B = 31406;   C = 41;    D = 30028;   E = 665;                 G = 2;   H = 713 }
    76  = @{ B = 19964;   C = 171;   D = 13095;   E = 6260;                G = 11;  H = 609 }
    77  = @{ B = 19947;   C = 248;   D = 14973;   E = 4650;                G = 1;   H = 324 }
    90  = @{ B = 10746;   C = 103;                E = 1134 }
    93  = @{ B = 10097;   C = 52;    D = 9133;    E = 706 }
    98  = @{ B = 8583;    C = 33;    D = 7376;    E = 1139 }
    109 = @{ B = 5387;    C = 2;     D = 5323;    E = 4 }
    184 = @{ B = 288;     C = 3;     D = 232;     E = 56 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# Timestamp bump
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 16:09"
